$p = $ppt.ActivePresentation

# The deck drops the "1 / Python Modules" section-divider slide (slide 2,
# which used the now-unused Subtitle/TITLE_1 layout). Every slide after it
# shifts up one position, so the slide-number field shown on each of those
# slides needs to be recomputed to match its new position.

$p.Slides.Item(2).Delete()

# Force PowerPoint to recompute the cached "slide number" field text on the
# remaining slides by toggling the footer's slide-number visibility off and
# back on - this regenerates the field with the slide's current number.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $hf = $slide.HeadersFooters
    if ($hf.SlideNumber.Visible) {
        $hf.SlideNumber.Visible = $false
        $hf.SlideNumber.Visible = $true
    }
}
